$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Row 6 key was renamed from "public.login.title" to "public.sign-in.title"
# (the translation text in C6 stays "Prihlaseni")
$ws.Cells.Item(6, 2).Value = "public.sign-in.title"

# Append 20 new translation rows (18-37) for client registration and
# related public pages. Copy formatting (wrap-text style) from the last
# existing data row so the new rows match the table look (style index).
$ws.Range("A17:C17").Copy()
$ws.Range("A18:C37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$rows = @(
  @("cs", "public.sign-up.menu", "Registrace"),
  @("cs", "public.404.title", "Tady nic není!"),
  @("cs", "public.404.title", "I těm nejlepším se to stane…"),
  @("cs", "public.404.subtitle", "Bohužel jste narazili na stránku, která neexistuje. Je to divné, ale je to tak."),
  @("cs", "public.404.back", "Zpět"),
  @("cs", "public.404.home", "Domů"),
  @("cs", "public.development-notice.alert", "Aplikace je stále ve vývoji a mnoho vychytávek chybí, nicméně je aktivně vyvíjena i používána."),
  @("cs", "public.sign-up.title", "Registrace"),
  @("cs", "public.sign-up.subtitle", "Po registraci získáte okamžitý přístup do aplikace a můžete začít zkoumat zajímavé možnosti, kterými disponuje."),
  @("cs", "user.name.label", "Vaše jméno"),
  @("cs", "user.name.label.tooltip", "Jedná se víceméně o jakékoli jméno, kterým se chcete prezentovat. Pro přihlášení bude použit Váš email."),
  @("cs", "user.password2.label", "Kontrola hesla"),
  @("cs", "public.sign-up.form.submit.label", "Registrovat se"),
  @("cs", "user.email.label", "E-mail"),
  @("cs", "user.email.label.tooltip", "Emailová adresa musí být v systému unikátní a slouží dále pro přihlášení do systému."),
  @("cs", "user.password.label.required", "Bez hesla to bohužel nejde, vyplňte jej prosím."),
  @("cs", "user.password2.label.required", "Kontrolní heslo je vyžadováno."),
  @("cs", "user.name.label.required", "Vaše jméno je vyžadováno, vyplňte jej prosím."),
  @("cs", "user.email.label.required", "E-mailová adresa slouží k přihlášení do systému, je povinná."),
  @("cs", "user.password2.label.mismatch", "Hesla se neshoduji!")
)

$r = 18
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# Match the selection recorded in the committed workbook
$ws.Range("B32").Select()
